$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 54 with revised quarterly values ---
$ws.Range("B54").Value = 81.3
$ws.Range("C54").Value = 81.8
$ws.Range("D54").Value = 98.59999999999999
$ws.Range("E54").Value = 79
$ws.Range("F54").Value = 91.90000000000001
$ws.Range("G54").Value = 93.90000000000001
$ws.Range("H54").Value = 133.4
$ws.Range("I54").Value = 135.3
$ws.Range("J54").Value = 145.7
$ws.Range("K54").Value = 114.3
$ws.Range("L54").Value = 126.3
$ws.Range("M54").Value = 116.2
$ws.Range("P54").Value = 108.7
$ws.Range("Q54").Value = 90.59999999999999
$ws.Range("R54").Value = 95.5
$ws.Range("S54").Value = 98.40000000000001
$ws.Range("U54").Value = 100.4
$ws.Range("X54").Value = 98.5
$ws.Range("Y54").Value = 103.1
$ws.Range("Z54").Value = 104.1
$ws.Range("AA54").Value = 102.1
$ws.Range("AC54").Value = 83.5
$ws.Range("AD54").Value = 83.3
$ws.Range("AE54").Value = 83.59999999999999
$ws.Range("AF54").Value = 82.90000000000001
$ws.Range("AG54").Value = 101.3
$ws.Range("AH54").Value = 77.40000000000001
$ws.Range("AJ54").Value = 103.2
$ws.Range("AL54").Value = 79.40000000000001
$ws.Range("AM54").Value = 109.1

# --- Insert new row 55 (01-04-2021) ---
$ws.Range("BZ100").Formula = "=""01-04-2021"""
$ws.Range("BZ100").Copy()
$ws.Range("A55").PasteSpecial(-4163)
$ws.Range("BZ100").ClearContents()

$ws.Range("B55").Value = 90
$ws.Range("C55").Value = 86.7
$ws.Range("D55").Value = 77.3
$ws.Range("E55").Value = 69.8
$ws.Range("F55").Value = 90.40000000000001
$ws.Range("G55").Value = 106.1
$ws.Range("H55").Value = 137
$ws.Range("I55").Value = 138.1
$ws.Range("J55").Value = 144.4
$ws.Range("K55").Value = 125.5
$ws.Range("L55").Value = 134.4
$ws.Range("M55").Value = 127
$ws.Range("N55").Value = 158.6
$ws.Range("O55").Value = 126.5
$ws.Range("P55").Value = 110.8
$ws.Range("Q55").Value = 95.7
$ws.Range("R55").Value = 103.2
$ws.Range("S55").Value = 112.2
$ws.Range("T55").Value = 91.8
$ws.Range("U55").Value = 103.7
$ws.Range("V55").Value = 97.90000000000001
$ws.Range("W55").Value = 97.90000000000001
$ws.Range("X55").Value = 96.59999999999999
$ws.Range("Y55").Value = 110.1
$ws.Range("Z55").Value = 111.1
$ws.Range("AA55").Value = 108.7
$ws.Range("AB55").Value = 105
$ws.Range("AC55").Value = 94.40000000000001
$ws.Range("AD55").Value = 95.90000000000001
$ws.Range("AE55").Value = 86
$ws.Range("AF55").Value = 85.7
$ws.Range("AG55").Value = 112.8
$ws.Range("AH55").Value = 77.59999999999999
$ws.Range("AI55").Value = 78.7
$ws.Range("AJ55").Value = 107.1
$ws.Range("AK55").Value = 95.2
$ws.Range("AL55").Value = 75.3
$ws.Range("AM55").Value = 117.4
